$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime the shared-string table so the two brand-new cluster names are interned
# in the same relative order the source data uses (ECs before M2); the four
# pre-existing strings (FAPs/sCs/Bmp2/Bmpr1b) keep the slots they already own.
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(2, 1).Value = "M2"

# Re-populate rows 2-9: one row per (sending cluster, target cluster) pair for the
# Bmp2 -> Bmpr1b ligand-receptor edge (NATMI lrc2p output, 4 sending clusters x 2 target clusters)
# Row 2: M2 -> FAPs
$ws.Cells.Item(2, 1).Value = "M2"
$ws.Cells.Item(2, 2).Value = "Bmp2"
$ws.Cells.Item(2, 3).Value = "Bmpr1b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.55635
$ws.Cells.Item(2, 8).Value = 1.66905
$ws.Cells.Item(2, 9).Value = 0.1200566932586554
$ws.Cells.Item(2, 10).Value = 0.1200566932586554
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.510190333333334
$ws.Cells.Item(2, 14).Value = 10.530571
$ws.Cells.Item(2, 15).Value = 0.8478537661184122
$ws.Cells.Item(2, 16).Value = 0.8478537661184122
$ws.Cells.Item(2, 17).Value = 1.95289439195
$ws.Cells.Item(2, 18).Value = 17.57604952755
$ws.Cells.Item(2, 19).Value = 0.1017905195270739
$ws.Cells.Item(2, 20).Value = 0.1017905195270739

# Row 3: M2 -> sCs
$ws.Cells.Item(3, 1).Value = "M2"
$ws.Cells.Item(3, 2).Value = "Bmp2"
$ws.Cells.Item(3, 3).Value = "Bmpr1b"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.55635
$ws.Cells.Item(3, 8).Value = 1.66905
$ws.Cells.Item(3, 9).Value = 0.1200566932586554
$ws.Cells.Item(3, 10).Value = 0.1200566932586554
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.629899
$ws.Cells.Item(3, 14).Value = 1.889697
$ws.Cells.Item(3, 15).Value = 0.1521462338815877
$ws.Cells.Item(3, 16).Value = 0.1521462338815877
$ws.Cells.Item(3, 17).Value = 0.35044430865
$ws.Cells.Item(3, 18).Value = 3.15399877785
$ws.Cells.Item(3, 19).Value = 0.01826617373158141
$ws.Cells.Item(3, 20).Value = 0.01826617373158142

# Row 4: sCs -> FAPs
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Bmp2"
$ws.Cells.Item(4, 3).Value = "Bmpr1b"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.611936333333333
$ws.Cells.Item(4, 8).Value = 4.835809
$ws.Cells.Item(4, 9).Value = 0.3478453238491627
$ws.Cells.Item(4, 10).Value = 0.3478453238491627
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.510190333333334
$ws.Cells.Item(4, 14).Value = 10.530571
$ws.Cells.Item(4, 15).Value = 0.8478537661184122
$ws.Cells.Item(4, 16).Value = 0.8478537661184122
$ws.Cells.Item(4, 17).Value = 5.658203335215445
$ws.Cells.Item(4, 18).Value = 50.923830016939
$ws.Cells.Item(4, 19).Value = 0.2949219678521913
$ws.Cells.Item(4, 20).Value = 0.2949219678521913

# Row 5: sCs -> sCs
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Bmp2"
$ws.Cells.Item(5, 3).Value = "Bmpr1b"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.611936333333333
$ws.Cells.Item(5, 8).Value = 4.835809
$ws.Cells.Item(5, 9).Value = 0.3478453238491627
$ws.Cells.Item(5, 10).Value = 0.3478453238491627
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.629899
$ws.Cells.Item(5, 14).Value = 1.889697
$ws.Cells.Item(5, 15).Value = 0.1521462338815877
$ws.Cells.Item(5, 16).Value = 0.1521462338815877
$ws.Cells.Item(5, 17).Value = 1.015357084430333
$ws.Cells.Item(5, 18).Value = 9.138213759873
$ws.Cells.Item(5, 19).Value = 0.05292335599697132
$ws.Cells.Item(5, 20).Value = 0.05292335599697133

# Row 6: ECs -> FAPs
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Bmp2"
$ws.Cells.Item(6, 3).Value = "Bmpr1b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.3793803333333334
$ws.Cells.Item(6, 8).Value = 1.138141
$ws.Cells.Item(6, 9).Value = 0.08186779600497247
$ws.Cells.Item(6, 10).Value = 0.08186779600497246
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.510190333333334
$ws.Cells.Item(6, 14).Value = 10.530571
$ws.Cells.Item(6, 15).Value = 0.8478537661184122
$ws.Cells.Item(6, 16).Value = 0.8478537661184122
$ws.Cells.Item(6, 17).Value = 1.331697178723445
$ws.Cells.Item(6, 18).Value = 11.985274608511
$ws.Cells.Item(6, 19).Value = 0.06941191916662981
$ws.Cells.Item(6, 20).Value = 0.06941191916662981

# Row 7: ECs -> sCs
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Bmp2"
$ws.Cells.Item(7, 3).Value = "Bmpr1b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.3793803333333334
$ws.Cells.Item(7, 8).Value = 1.138141
$ws.Cells.Item(7, 9).Value = 0.08186779600497247
$ws.Cells.Item(7, 10).Value = 0.08186779600497246
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.629899
$ws.Cells.Item(7, 14).Value = 1.889697
$ws.Cells.Item(7, 15).Value = 0.1521462338815877
$ws.Cells.Item(7, 16).Value = 0.1521462338815877
$ws.Cells.Item(7, 17).Value = 0.2389712925863333
$ws.Cells.Item(7, 18).Value = 2.150741633277
$ws.Cells.Item(7, 19).Value = 0.01245587683834265
$ws.Cells.Item(7, 20).Value = 0.01245587683834265

# Row 8: FAPs -> FAPs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Bmp2"
$ws.Cells.Item(8, 3).Value = "Bmpr1b"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.086394
$ws.Cells.Item(8, 8).Value = 6.259182
$ws.Cells.Item(8, 9).Value = 0.4502301868872095
$ws.Cells.Item(8, 10).Value = 0.4502301868872095
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.510190333333334
$ws.Cells.Item(8, 14).Value = 10.530571
$ws.Cells.Item(8, 15).Value = 0.8478537661184122
$ws.Cells.Item(8, 16).Value = 0.8478537661184122
$ws.Cells.Item(8, 17).Value = 7.323640050324666
$ws.Cells.Item(8, 18).Value = 65.912760452922
$ws.Cells.Item(8, 19).Value = 0.3817293595725171
$ws.Cells.Item(8, 20).Value = 0.3817293595725171

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Bmp2"
$ws.Cells.Item(9, 3).Value = "Bmpr1b"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.086394
$ws.Cells.Item(9, 8).Value = 6.259182
$ws.Cells.Item(9, 9).Value = 0.4502301868872095
$ws.Cells.Item(9, 10).Value = 0.4502301868872095
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.629899
$ws.Cells.Item(9, 14).Value = 1.889697
$ws.Cells.Item(9, 15).Value = 0.1521462338815877
$ws.Cells.Item(9, 16).Value = 0.1521462338815877
$ws.Cells.Item(9, 17).Value = 1.314217494206
$ws.Cells.Item(9, 18).Value = 11.827957447854
$ws.Cells.Item(9, 19).Value = 0.06850082731469233
$ws.Cells.Item(9, 20).Value = 0.06850082731469233
